$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the H1 title
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleEnd = $titlePara.Range
$titleEnd.Collapse(0)
$titleEnd.InsertParagraphAfter()

# The freshly created paragraph inherited the Heading1 style - reset it back
# to the plain body style used elsewhere in the document.
$d.Paragraphs.Item(2).Style = $d.Styles.Item("Normal")

$newParaStart = $d.Paragraphs.Item(2).Range.Start
$insertionPoint = $d.Range($newParaStart, $newParaStart)

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our unbiased review of Big Bucks Bandits Megaways slot game. Get ratings, betting options, and play for free!</w:t></w:r></w:p>' + `
    '</w:body></w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2. Remove the duplicated bold title paragraph near the end of the document.
#    Locate it by matching the title text again (rather than assuming a
#    fixed index) so the script stays correct even if content shifts.
# ---------------------------------------------------------------------------
$titleText = $d.Paragraphs.Item(1).Range.Text
$dupIndex = -1
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq $titleText) {
        $dupIndex = $i
    }
}

if ($dupIndex -gt 0) {
    $dupTitlePara = $d.Paragraphs.Item($dupIndex)
    $dupTitleRange = $d.Range($dupTitlePara.Range.Start, $dupTitlePara.Range.End)
    $dupTitleRange.Delete()
}

# ---------------------------------------------------------------------------
# 3. Replace the old meta-description sentence (now an image prompt) while
#    keeping its italic formatting intact. Scope the search to the very last
#    paragraph only, since the sentence also now appears earlier in the
#    document (in the newly inserted "Meta description" paragraph).
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$find = $lastPara.Range.Find
$find.Execute(
    "Read our unbiased review of Big Bucks Bandits Megaways slot game. Get ratings, betting options, and play for free!",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Create an exciting feature image for Big Bucks Bandits Megaways that will catch the attention of online slot players. The image should be in a cartoon style featuring a happy Maya warrior with glasses. The warrior should be dressed in traditional clothing, with a headdress adorned with gold. In the background, there should be a desert landscape with the Grand Canyon visible. The warrior should be holding a winning slot combination of symbols (such as a horse, guns, cacti, and the like) with an excited expression on their face. The image should be vibrant and colorful to capture the excitement and spirit of the Wild West. This image will entice players to try their luck at Big Bucks Bandits Megaways and experience the thrill of hitting the jackpot.",
    2)
